# Burndown Chart sprint-2 update: shift some daily task-effort entries to
# different days/users and bump a couple totals, per the "Atualiza lógica de
# inserção de usuário e ajusta visibilidade de botões no layout" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")
$ws.Activate() | Out-Null

# --- Row 2 (S1 - Daniel): the "1" on Dia 9 (K) moves to Dia 8 (J); Dia 11 (M) gets a new 3
$ws.Range("J2").Value = 1
$ws.Range("K2").ClearContents() | Out-Null
$ws.Range("M2").Value = 3

# --- Row 3 (S2.1 - Front): the "1" on Dia 8 (J) moves to Dia 7 (I)
$ws.Range("I3").Value = 1
$ws.Range("J3").ClearContents() | Out-Null

# --- Row 4 (S2.2 - Front): Dia 11 (M) now logs 1.5h
$ws.Range("M4").Value = 1.5

# --- Row 6 (S3.1 - Front): Dia 5 (G)->Dia 4 (F); Dia 8 (J)->Dia 7 (I)
$ws.Range("F6").Value = 1
$ws.Range("G6").ClearContents() | Out-Null
$ws.Range("I6").Value = 4
$ws.Range("J6").ClearContents() | Out-Null

# --- Row 10 (S5 - Bueno): Dia 9 (K)->Dia 8 (J)
$ws.Range("J10").Value = 4
$ws.Range("K10").ClearContents() | Out-Null

# --- Row 13 (S8.1 - Rafael): Dia 4 (F)->Dia 3 (E)
$ws.Range("E13").Value = 1
$ws.Range("F13").ClearContents() | Out-Null

# --- Row 14 (S8.2 - Rafael): Dia 4 (F)->Dia 3 (E)
$ws.Range("E14").Value = 1
$ws.Range("F14").ClearContents() | Out-Null

# --- Row 15 (S8.3 - Rafael): everything shifts one day earlier (E gets the old F value,
# F gets the old G value, G is cleared)
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("G15").ClearContents() | Out-Null

# --- Row 16 (S9 - Rafael): total hours 8 -> 13, and logged hours move around
$ws.Range("B16").Value = 13
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 3
$ws.Range("I16").ClearContents() | Out-Null
$ws.Range("L16").Value = 2
$ws.Range("M16").Value = 4

# --- Row 23 (S2 - Estudo): total hours 8 -> 5, plus 1h logged on Dia 10 (L)
$ws.Range("B23").Value = 5
$ws.Range("L23").Value = 1

# --- Row 24 (S0 - Git e docs): the 1h on Dia 11 (M) moves to Dia 10 (L)
$ws.Range("L24").Value = 1
$ws.Range("M24").ClearContents() | Out-Null

# --- Row 26 (Estimado): B26 is a literal seed value (not a formula); the rest of the
# row recalculates off of it automatically.
$ws.Range("B26").Value = 106

$wb.Application.Calculate() | Out-Null

# Restore the author's final selection on the sheet
$ws.Range("M9").Select() | Out-Null
